# Updated symbol list on Fri Jan 20 08:49:51 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns of the crypto table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'288.35"
$ws.Range('E2').Value = "'-1.03%"
$ws.Range('D3').Value = "'31.03"
$ws.Range('E3').Value = "'1.25%"
$ws.Range('E4').Value = "'-0.51%"
$ws.Range('E5').Value = "'1.70%"
$ws.Range('D6').Value = "'2.207"
$ws.Range('E6').Value = "'22.91%"
$ws.Range('D7').Value = "'7.717"
$ws.Range('E7').Value = "'0.39%"
$ws.Range('D8').Value = "'3.733"
$ws.Range('E8').Value = "'-0.62%"
$ws.Range('D9').Value = "'0.9039"
$ws.Range('D10').Value = "'0.09188"
$ws.Range('E10').Value = "'19.21%"
$ws.Range('D11').Value = "'0.1680"
$ws.Range('E11').Value = "'1.15%"
$ws.Range('D12').Value = "'0.08260"
$ws.Range('E12').Value = "'2.02%"
$ws.Range('E13').Value = "'2.41%"
$ws.Range('D14').Value = "'0.09959"
$ws.Range('E14').Value = "'-0.59%"
$ws.Range('E15').Value = "'0.07%"
$ws.Range('D16').Value = "'0.005754"
$ws.Range('E16').Value = "'1.61%"
$ws.Range('D17').Value = "'3.521"
$ws.Range('E17').Value = "'1.34%"
$ws.Range('E19').Value = "'0.48%"
$ws.Range('E20').Value = "'0.40%"
$ws.Range('D21').Value = "'4.168"
$ws.Range('E21').Value = "'3.02%"
$ws.Range('D22').Value = "'0.2102"
$ws.Range('E22').Value = "'-6.68%"
$ws.Range('D23').Value = "'0.04542"
$ws.Range('E23').Value = "'0.94%"
$ws.Range('E24').Value = "'-0.59%"
$ws.Range('D25').Value = "'0.004158"
$ws.Range('E25').Value = "'3.63%"
$ws.Range('E26').Value = "'3.98%"
$ws.Range('D27').Value = "'0.0003398"
$ws.Range('D39').Value = "'0.01568"
$ws.Range('E39').Value = "'-2.46%"
$ws.Range('D40').Value = "'0.04438"
$ws.Range('E40').Value = "'0.39%"
$ws.Range('D41').Value = "'0.007367"
$ws.Range('E41').Value = "'1.13%"
$ws.Range('D42').Value = "'0.008998"
$ws.Range('E42').Value = "'3.31%"
$ws.Range('D43').Value = "'0.1329"
$ws.Range('E43').Value = "'1.58%"
$ws.Range('D44').Value = "'0.002232"
$ws.Range('E44').Value = "'10.75%"
$ws.Range('E45').Value = "'-4.37%"
$ws.Range('D46').Value = "'0.00006117"
$ws.Range('E46').Value = "'2.77%"
$ws.Range('E47').Value = "'-0.01%"
$ws.Range('D48').Value = "'2.162"
$ws.Range('E48').Value = "'-3.78%"
$ws.Range('D49').Value = "'0.002002"
$ws.Range('E49').Value = "'-33.34%"
$ws.Range('D50').Value = "'0.00002102"
$ws.Range('E50').Value = "'-0.01%"
$ws.Range('D51').Value = "'0.0002002"
$ws.Range('E51').Value = "'-0.01%"
